# Auto update Excel log
# Appends newly-logged sensor readings (2026-01-28 afternoon batch) to the
# PIR, Humidity, Temperature and Proximity sheets of the SeniorConnect
# master log.

$wb = $excel.ActiveWorkbook

# --- Pre-format the columns that hold "date-looking" / "percent-looking"
# text so Excel stores them as literal text instead of silently coercing
# them into date serials / percentages (matches how the rest of the log
# is already stored: plain text cells).
$wsPIR = $wb.Worksheets.Item("PIR")
$wsPIR.Range("A258:A270").NumberFormat = "@"

$wsHumidity = $wb.Worksheets.Item("Humidity")
$wsHumidity.Range("A248:A258").NumberFormat = "@"
$wsHumidity.Range("E248:E258").NumberFormat = "@"

$wsTemperature = $wb.Worksheets.Item("Temperature")
$wsTemperature.Range("A249:A259").NumberFormat = "@"

$wsProximity = $wb.Worksheets.Item("Proximity")
$wsProximity.Range("A2:A2").NumberFormat = "@"

# --- PIR sheet: rows 258-270, Bathroom / No Motion / Inactive -------------
$pirTimes = @(
    "16:31:39","16:31:40","16:31:46","16:31:50","16:31:55","16:32:00",
    "16:32:05","16:32:10","16:32:15","16:32:20","16:32:26","16:32:30",
    "16:32:35"
)
$r = 258
foreach ($t in $pirTimes) {
    $wsPIR.Cells.Item($r, 1).Value = "2026-01-28"
    $wsPIR.Cells.Item($r, 2).Value = $t
    $wsPIR.Cells.Item($r, 3).Value = "16:00"
    $wsPIR.Cells.Item($r, 4).Value = "Bathroom"
    $wsPIR.Cells.Item($r, 5).Value = "No Motion"
    $wsPIR.Cells.Item($r, 6).Value = "Inactive"
    $r = $r + 1
}

# --- Humidity sheet: rows 248-258, Bathroom / xx.x% / Active --------------
$humidityRows = @(
    @("16:31:40","88.0%"),
    @("16:31:44","87.1%"),
    @("16:31:48","88.1%"),
    @("16:31:52","88.0%"),
    @("16:31:56","88.0%"),
    @("16:32:08","88.0%"),
    @("16:32:12","88.0%"),
    @("16:32:21","87.9%"),
    @("16:32:24","87.0%"),
    @("16:32:32","87.9%"),
    @("16:32:36","87.0%")
)
$r = 248
foreach ($row in $humidityRows) {
    $wsHumidity.Cells.Item($r, 1).Value = "2026-01-28"
    $wsHumidity.Cells.Item($r, 2).Value = $row[0]
    $wsHumidity.Cells.Item($r, 3).Value = "16:00"
    $wsHumidity.Cells.Item($r, 4).Value = "Bathroom"
    $wsHumidity.Cells.Item($r, 5).Value = $row[1]
    $wsHumidity.Cells.Item($r, 6).Value = "Active"
    $r = $r + 1
}

# --- Temperature sheet: rows 249-259, Bathroom / xx.xC / Active -----------
$temperatureRows = @(
    @("16:31:41","22.9C"),
    @("16:31:45","22.9C"),
    @("16:31:49","22.8C"),
    @("16:31:53","22.9C"),
    @("16:31:57","22.8C"),
    @("16:32:09","22.8C"),
    @("16:32:13","22.8C"),
    @("16:32:21","22.8C"),
    @("16:32:25","22.8C"),
    @("16:32:33","22.8C"),
    @("16:32:37","22.8C")
)
$r = 249
foreach ($row in $temperatureRows) {
    $wsTemperature.Cells.Item($r, 1).Value = "2026-01-28"
    $wsTemperature.Cells.Item($r, 2).Value = $row[0]
    $wsTemperature.Cells.Item($r, 3).Value = "16:00"
    $wsTemperature.Cells.Item($r, 4).Value = "Bathroom"
    $wsTemperature.Cells.Item($r, 5).Value = $row[1]
    $wsTemperature.Cells.Item($r, 6).Value = "Active"
    $r = $r + 1
}

# --- Proximity sheet: row 2, Living Room / Presence Detected / Active ----
$wsProximity.Cells.Item(2, 1).Value = "2026-01-28"
$wsProximity.Cells.Item(2, 2).Value = "16:32:34"
$wsProximity.Cells.Item(2, 3).Value = "16:00"
$wsProximity.Cells.Item(2, 4).Value = "Living Room"
$wsProximity.Cells.Item(2, 5).Value = "Presence Detected"
$wsProximity.Cells.Item(2, 6).Value = "Active"

Write-Host "SeniorConnect master log updated: PIR +13, Humidity +11, Temperature +11, Proximity +1"
